$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.357.10"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "1.938.55"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("D5").Value = "'250.58"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "'0.7257"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.3338"
$ws.Range("E8").Value = "  -4.63%  "
$ws.Range("D9").Value = "'28.47"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").Value = "'0.07252"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("D11").Value = "'0.8133"
$ws.Range("E11").Value = "  -3.89%  "
$ws.Range("D12").Value = "'0.08099"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "1.938.97"
$ws.Range("D14").Value = "'5.487"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "'94.55"
$ws.Range("E15").Value = "  -6.53%  "
$ws.Range("D16").Value = "'15.02"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "30.365.71"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "'0.000008265"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").Value = "'250.23"
$ws.Range("E19").Value = "  -8.26%  "
$ws.Range("D20").Value = "'5.914"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "2.193.61"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'6.969"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").Value = "'9.772"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").Value = "'163.40"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "'2.399"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'19.30"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "'0.1330"
$ws.Range("E29").Value = "  -8.06%  "
$ws.Range("D30").Value = "'1.571"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "'1.348"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "'4.453"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").Value = "'4.200"
$ws.Range("E33").Value = "  -5.34%  "
$ws.Range("D34").Value = "'0.05219"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "'1.296"
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("D36").Value = "'0.7511"
$ws.Range("E36").Value = "  -4.67%  "
$ws.Range("D37").Value = "'2.749"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'0.01986"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").Value = "'2.832"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").Value = "'81.00"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "'6.465"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").Value = "'0.4559"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "'2.048"
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").Value = "'0.8481"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "'102.17"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").Value = "'9.816"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "'7.466"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").Value = "'36.94"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'2.885"
$ws.Range("E50").Value = "  +7.00%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4207"
$ws.Range("E51").Value = "  -2.88%  "
